$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.175.01"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "3.740.11"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "3.737.71"
$ws.Range("E7").Value = "  -3.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.722"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -8.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000295"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.25%  "
$ws.Range("D15").Value = "4.342.16"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "3.746.85"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.127"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.76%  "
$ws.Range("D21").Value = "69.025.95"
$ws.Range("E21").Value = "  -3.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "415.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.32%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.59%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.76%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -18.00%  "
$ws.Range("E33").Value = "  -6.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.122"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "624.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "44.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.35%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "66.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("D38").Value = "0.0₃0893"
$ws.Range("E38").Value = "  -10.53%  "
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.405"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0444"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -17.98%  "
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.27%  "
$ws.Range("D49").Value = "2.824.89"
$ws.Range("E49").Value = "  -3.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.54%  "
